$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2005.1714
$ws.Range("I40").Value = 1922.2273
$ws.Range("J40").Value = 2145.5386
$ws.Range("K40").Value = 1922.2273
$ws.Range("L40").Value = 2145.5386
$ws.Range("M40").Value = -1747.2273
$ws.Range("N40").Value = -2495.5386

$ws.Range("H43").Value = 632.6
$ws.Range("I43").Value = 645.5
$ws.Range("K43").Value = 645.5
$ws.Range("M43").Value = -576.5

$ws.Range("H62").Value = 9886.388999999999
$ws.Range("I62").Value = 7843.077
$ws.Range("J62").Value = 15199
$ws.Range("K62").Value = 7843.077
$ws.Range("L62").Value = 15199
$ws.Range("M62").Value = -7219.077
$ws.Range("N62").Value = -16447

$ws.Range("H65").Value = 9886.388999999999
$ws.Range("I65").Value = 7843.077
$ws.Range("J65").Value = 15199
$ws.Range("K65").Value = 39215.385
$ws.Range("L65").Value = 75995
$ws.Range("M65").Value = -36095.385
$ws.Range("N65").Value = -82235

$ws.Range("H132").Value = 344322.22
$ws.Range("I132").Value = 450984.66
$ws.Range("J132").Value = 56333.6
$ws.Range("K132").Value = 1352953.98
$ws.Range("L132").Value = 169000.8
$ws.Range("M132").Value = -1350423.98
$ws.Range("N132").Value = -174060.8

$ws.Range("H133").Value = 27999.5
$ws.Range("J133").Value = 27999.5
$ws.Range("L133").Value = 27999.5
$ws.Range("N133").Value = -38119.5

$ws.Range("H137").Value = 26317118
$ws.Range("I137").Value = 32258968
$ws.Range("J137").Value = 3205.2856
$ws.Range("K137").Value = 96776904
$ws.Range("L137").Value = 9615.856800000001
$ws.Range("M137").Value = -96774354
$ws.Range("N137").Value = -14715.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 182.75
$ws.Range("I4").Value = 216.2
$ws.Range("J4").Value = 127
$ws.Range("K4").Value = 216.2
$ws.Range("L4").Value = 127
$ws.Range("M4").Value = -100.2
$ws.Range("N4").Value = -359

$ws.Range("H5").Value = 1000280.1
$ws.Range("J5").Value = 212.5
$ws.Range("L5").Value = 212.5
$ws.Range("N5").Value = -436.5

$ws.Range("H107").Value = 25001
$ws.Range("J107").Value = 25001
$ws.Range("L107").Value = 25001
$ws.Range("N107").Value = -32681

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1000280.1
$ws.Range("J4").Value = 212.5
$ws.Range("L4").Value = 212.5
$ws.Range("N4").Value = -442.5

$ws.Range("H94").Value = 968.8095
$ws.Range("I94").Value = 1071.7142
$ws.Range("J94").Value = 763
$ws.Range("K94").Value = 1071.7142
$ws.Range("L94").Value = 763
$ws.Range("M94").Value = -620.7141999999999
$ws.Range("N94").Value = -1665

$ws.Range("H99").Value = 1227.3
$ws.Range("I99").Value = 998.8333
$ws.Range("J99").Value = 1570
$ws.Range("K99").Value = 998.8333
$ws.Range("L99").Value = 1570
$ws.Range("M99").Value = 499.1667
$ws.Range("N99").Value = -4566

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1209.1143
$ws.Range("I31").Value = 947.78125
$ws.Range("K31").Value = 947.78125
$ws.Range("M31").Value = -652.78125

$ws.Range("H34").Value = 1209.1143
$ws.Range("I34").Value = 947.78125
$ws.Range("K34").Value = 947.78125
$ws.Range("M34").Value = -745.78125

$ws.Range("H99").Value = 7813794
$ws.Range("I99").Value = 12501040
$ws.Range("J99").Value = 1716.6666
$ws.Range("K99").Value = 12501040
$ws.Range("L99").Value = 1716.6666
$ws.Range("M99").Value = -12499542
$ws.Range("N99").Value = -4712.6666

$ws.Range("H126").Value = 7813794
$ws.Range("I126").Value = 12501040
$ws.Range("J126").Value = 1716.6666
$ws.Range("K126").Value = 37503120
$ws.Range("L126").Value = 5149.9998
$ws.Range("M126").Value = -37500650
$ws.Range("N126").Value = -10089.9998

$ws.Range("H132").Value = 2150.1538
$ws.Range("I132").Value = 1657.375
$ws.Range("J132").Value = 4402.857
$ws.Range("K132").Value = 4972.125
$ws.Range("L132").Value = 13208.571
$ws.Range("M132").Value = -2442.125
$ws.Range("N132").Value = -18268.571

$ws.Range("H134").Value = 2548.9678
$ws.Range("I134").Value = 1314.96
$ws.Range("J134").Value = 7690.6665
$ws.Range("K134").Value = 3944.88
$ws.Range("L134").Value = 23071.9995
$ws.Range("M134").Value = -1409.88
$ws.Range("N134").Value = -28141.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 11561.363
$ws.Range("I68").Value = 24464.4
$ws.Range("J68").Value = 808.8333
$ws.Range("K68").Value = 73393.20000000001
$ws.Range("L68").Value = 2426.4999
$ws.Range("M68").Value = -72582.20000000001
$ws.Range("N68").Value = -4048.4999

$ws.Range("H71").Value = 11561.363
$ws.Range("I71").Value = 24464.4
$ws.Range("J71").Value = 808.8333
$ws.Range("K71").Value = 220179.6
$ws.Range("L71").Value = 7279.4997
$ws.Range("M71").Value = -216123.6
$ws.Range("N71").Value = -15391.4997

$ws.Range("H139").Value = 2110.5667
$ws.Range("I139").Value = 1600.68
$ws.Range("J139").Value = 4660
$ws.Range("K139").Value = 4802.04
$ws.Range("L139").Value = 13980
$ws.Range("M139").Value = 337.96
$ws.Range("N139").Value = -24260

$ws.Range("H140").Value = 5006.0347
$ws.Range("I140").Value = 6704.4116
$ws.Range("K140").Value = 20113.2348
$ws.Range("M140").Value = -14933.2348

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2628.2856
$ws.Range("I80").Value = 2350
$ws.Range("J80").Value = 2739.6
$ws.Range("K80").Value = 2350
$ws.Range("L80").Value = 2739.6
$ws.Range("M80").Value = -1352
$ws.Range("N80").Value = -4735.6

$ws.Range("H83").Value = 2628.2856
$ws.Range("I83").Value = 2350
$ws.Range("J83").Value = 2739.6
$ws.Range("K83").Value = 11750
$ws.Range("L83").Value = 13698
$ws.Range("M83").Value = -6758
$ws.Range("N83").Value = -23682

$ws.Range("H123").Value = 10532.583
$ws.Range("J123").Value = 10532.583
$ws.Range("L123").Value = 10532.583
$ws.Range("N123").Value = -15432.583

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 424.2
$ws.Range("I55").Value = 335.875
$ws.Range("J55").Value = 525.1429000000001
$ws.Range("K55").Value = 335.875
$ws.Range("L55").Value = 525.1429000000001
$ws.Range("M55").Value = -162.875
$ws.Range("N55").Value = -871.1429000000001

$ws.Range("H122").Value = 3207.1428
$ws.Range("I122").Value = 2263.182
$ws.Range("J122").Value = 3817.9412
$ws.Range("K122").Value = 6789.545999999999
$ws.Range("L122").Value = 11453.8236
$ws.Range("M122").Value = -4339.545999999999
$ws.Range("N122").Value = -16353.8236

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 49900
$ws.Range("J108").Value = 49900
$ws.Range("L108").Value = 49900
$ws.Range("N108").Value = -57580
